$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 3900
$ws.Range("J40").Value = 3900
$ws.Range("L40").Value = 3900
$ws.Range("N40").Value = -4250

$ws.Range("H43").Value = 5010
$ws.Range("I43").Value = 5450
$ws.Range("K43").Value = 5450
$ws.Range("M43").Value = -5381

$ws.Range("H98").Value = 4559.5
$ws.Range("I98").Value = 4803.3076
$ws.Range("J98").Value = 1390
$ws.Range("K98").Value = 4803.3076
$ws.Range("L98").Value = 1390
$ws.Range("M98").Value = -3305.3076
$ws.Range("N98").Value = -4386

$ws.Range("H116").Value = 31262122
$ws.Range("I116").Value = 125004990
$ws.Range("J116").Value = 14498.333
$ws.Range("K116").Value = 125004990
$ws.Range("L116").Value = 14498.333
$ws.Range("M116").Value = -125001548
$ws.Range("N116").Value = -21382.333

$ws.Range("H122").Value = 4559.5
$ws.Range("I122").Value = 4803.3076
$ws.Range("J122").Value = 1390
$ws.Range("K122").Value = 14409.9228
$ws.Range("L122").Value = 1390
$ws.Range("M122").Value = -11959.9228
$ws.Range("N122").Value = -9070

$ws.Range("H125").Value = 100001620
$ws.Range("I125").Value = 166667790
$ws.Range("J125").Value = 2361.75
$ws.Range("K125").Value = 1500010110
$ws.Range("L125").Value = 21255.75
$ws.Range("M125").Value = -1500007650
$ws.Range("N125").Value = -26175.75

$ws.Range("H127").Value = 711.5714
$ws.Range("I127").Value = 711.5714
$ws.Range("K127").Value = 2134.7142
$ws.Range("M127").Value = 2825.2858

$ws.Range("H129").Value = 1407.4348
$ws.Range("I129").Value = 895.6
$ws.Range("J129").Value = 1801.1538
$ws.Range("K129").Value = 2686.8
$ws.Range("L129").Value = 5403.4614
$ws.Range("M129").Value = 2313.2
$ws.Range("N129").Value = -15403.4614

$ws.Range("H137").Value = 5683.926
$ws.Range("I137").Value = 3205.6
$ws.Range("J137").Value = 8781.833000000001
$ws.Range("K137").Value = 9616.799999999999
$ws.Range("L137").Value = 26345.499
$ws.Range("M137").Value = -7066.799999999999
$ws.Range("N137").Value = -31445.499

$ws.Range("H138").Value = 1698972.5
$ws.Range("I138").Value = 2220.111
$ws.Range("K138").Value = 6660.333
$ws.Range("M138").Value = -1520.333

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 4262130.5
$ws.Range("I32").Value = 4550332.5
$ws.Range("J32").Value = 35166.332
$ws.Range("K32").Value = 4550332.5
$ws.Range("L32").Value = 35166.332
$ws.Range("M32").Value = -4550045.5
$ws.Range("N32").Value = -35740.332

$ws.Range("H36").Value = 12000
$ws.Range("I36").Value = 12000
$ws.Range("K36").Value = 12000
$ws.Range("M36").Value = -11654

$ws.Range("H61").Value = 10339.934
$ws.Range("I61").Value = 2418.8
$ws.Range("J61").Value = 14300.5
$ws.Range("K61").Value = 2418.8
$ws.Range("L61").Value = 14300.5
$ws.Range("M61").Value = -2206.8
$ws.Range("N61").Value = -14724.5

$ws.Range("H74").Value = 34983.453
$ws.Range("I74").Value = 49460.668
$ws.Range("J74").Value = 4581.3
$ws.Range("K74").Value = 49460.668
$ws.Range("L74").Value = 4581.3
$ws.Range("M74").Value = -48586.668
$ws.Range("N74").Value = -6329.3

$ws.Range("H77").Value = 34983.453
$ws.Range("I77").Value = 49460.668
$ws.Range("J77").Value = 4581.3
$ws.Range("K77").Value = 247303.34
$ws.Range("L77").Value = 22906.5
$ws.Range("M77").Value = -242935.34
$ws.Range("N77").Value = -31642.5

$ws.Range("H102").Value = 3442.625
$ws.Range("I102").Value = 2801.6
$ws.Range("K102").Value = 2801.6
$ws.Range("M102").Value = -1179.6

$ws.Range("H132").Value = 6609.5283
$ws.Range("I132").Value = 5044.3335
$ws.Range("J132").Value = 9924.058999999999
$ws.Range("K132").Value = 15133.0005
$ws.Range("L132").Value = 29772.177
$ws.Range("M132").Value = -12603.0005
$ws.Range("N132").Value = -34832.177

$ws.Range("H136").Value = 10339.934
$ws.Range("I136").Value = 2418.8
$ws.Range("J136").Value = 14300.5
$ws.Range("K136").Value = 7256.400000000001
$ws.Range("L136").Value = 42901.5
$ws.Range("M136").Value = -4706.400000000001
$ws.Range("N136").Value = -48001.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 41672524
$ws.Range("I20").Value = 53033490
$ws.Range("K20").Value = 53033490
$ws.Range("M20").Value = -53033243

$ws.Range("H105").Value = 6050.075
$ws.Range("I105").Value = 6666.375
$ws.Range("J105").Value = 5125.625
$ws.Range("K105").Value = 6666.375
$ws.Range("L105").Value = 5125.625
$ws.Range("M105").Value = -4919.375
$ws.Range("N105").Value = -8619.625

$ws.Range("H107").Value = 66183736
$ws.Range("I107").Value = 75007816
$ws.Range("J107").Value = 3149.5
$ws.Range("K107").Value = 75007816
$ws.Range("L107").Value = 3149.5
$ws.Range("M107").Value = -75005896
$ws.Range("N107").Value = -6989.5

$ws.Range("H134").Value = 4364.849
$ws.Range("I134").Value = 1622.1578
$ws.Range("K134").Value = 4866.4734
$ws.Range("M134").Value = -2331.4734

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 6433.0586
$ws.Range("I31").Value = 1813.3334
$ws.Range("J31").Value = 10539.481
$ws.Range("K31").Value = 1813.3334
$ws.Range("L31").Value = 10539.481
$ws.Range("M31").Value = -1518.3334
$ws.Range("N31").Value = -11129.481

$ws.Range("H34").Value = 6433.0586
$ws.Range("I34").Value = 1813.3334
$ws.Range("J34").Value = 10539.481
$ws.Range("K34").Value = 1813.3334
$ws.Range("L34").Value = 10539.481
$ws.Range("M34").Value = -1611.3334
$ws.Range("N34").Value = -10943.481

$ws.Range("H58").Value = 5308.1963
$ws.Range("I58").Value = 1946.76
$ws.Range("J58").Value = 8540.346
$ws.Range("K58").Value = 1946.76
$ws.Range("L58").Value = 8540.346
$ws.Range("M58").Value = -1743.76
$ws.Range("N58").Value = -8946.346

$ws.Range("H132").Value = 8903.190000000001
$ws.Range("I132").Value = 6138.857
$ws.Range("K132").Value = 18416.571
$ws.Range("M132").Value = -15886.571

$ws.Range("H134").Value = 8795.666999999999
$ws.Range("I134").Value = 2761.4
$ws.Range("K134").Value = 8284.200000000001
$ws.Range("M134").Value = -5749.200000000001

$ws.Range("H136").Value = 5308.1963
$ws.Range("I136").Value = 1946.76
$ws.Range("J136").Value = 8540.346
$ws.Range("K136").Value = 5840.28
$ws.Range("L136").Value = 25621.038
$ws.Range("M136").Value = -3290.28
$ws.Range("N136").Value = -30721.038

$ws.Range("H141").Value = 403599.6
$ws.Range("J141").Value = 403599.6
$ws.Range("L141").Value = 403599.6
$ws.Range("N141").Value = -413959.6

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 30827320
$ws.Range("I4").Value = 45837708
$ws.Range("K4").Value = 137513124
$ws.Range("M4").Value = -137513012

$ws.Range("H18").Value = 881.4286
$ws.Range("I18").Value = 421
$ws.Range("J18").Value = 2032.5
$ws.Range("K18").Value = 1263
$ws.Range("L18").Value = 6097.5
$ws.Range("M18").Value = -1094
$ws.Range("N18").Value = -6435.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 58830900
$ws.Range("I70").Value = 83339160
$ws.Range("K70").Value = 83339160
$ws.Range("M70").Value = -83338890

$ws.Range("H73").Value = 58830900
$ws.Range("I73").Value = 83339160
$ws.Range("K73").Value = 83339160
$ws.Range("M73").Value = -83338224

$ws.Range("H97").Value = 1655.1666
$ws.Range("I97").Value = 1739.4
$ws.Range("K97").Value = 1739.4
$ws.Range("M97").Value = -1243.4

$ws.Range("H102").Value = 3789.8823
$ws.Range("I102").Value = 3602.1724
$ws.Range("J102").Value = 4878.6
$ws.Range("K102").Value = 3602.1724
$ws.Range("L102").Value = 4878.6
$ws.Range("M102").Value = -1980.1724
$ws.Range("N102").Value = -8122.6

$ws.Range("H113").Value = 6854.1113
$ws.Range("I113").Value = 3750.611
$ws.Range("K113").Value = 3750.611
$ws.Range("M113").Value = -1580.611

$ws.Range("H132").Value = 5383.4
$ws.Range("I132").Value = 2662.4375
$ws.Range("K132").Value = 7987.3125
$ws.Range("M132").Value = -5457.3125

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H7").Value = 4811.467
$ws.Range("I7").Value = 2719.7778
$ws.Range("J7").Value = 7949
$ws.Range("K7").Value = 2719.7778
$ws.Range("L7").Value = 7949
$ws.Range("M7").Value = -2607.7778
$ws.Range("N7").Value = -8173

$ws.Range("H68").Value = 5047.3335
$ws.Range("J68").Value = 6889.9
$ws.Range("L68").Value = 6889.9
$ws.Range("N68").Value = -8387.9

$ws.Range("H71").Value = 5047.3335
$ws.Range("J71").Value = 6889.9
$ws.Range("L71").Value = 34449.5
$ws.Range("N71").Value = -41937.5

$ws.Range("H93").Value = 749.875
$ws.Range("I93").Value = 642.7143
$ws.Range("K93").Value = 642.7143
$ws.Range("M93").Value = 605.2857

$ws.Range("H126").Value = 4811.467
$ws.Range("I126").Value = 2719.7778
$ws.Range("J126").Value = 7949
$ws.Range("K126").Value = 8159.3334
$ws.Range("L126").Value = 23847
$ws.Range("M126").Value = -5689.3334
$ws.Range("N126").Value = -28787

$ws.Range("H132").Value = 6343.4346
$ws.Range("I132").Value = 3299.923
$ws.Range("J132").Value = 10300
$ws.Range("K132").Value = 9899.769
$ws.Range("L132").Value = 30900
$ws.Range("M132").Value = -7369.769
$ws.Range("N132").Value = -35960

$ws.Range("H136").Value = 11552.898
$ws.Range("I136").Value = 2943.96
$ws.Range("K136").Value = 8831.880000000001
$ws.Range("M136").Value = -6281.880000000001

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 737.3333
$ws.Range("I100").Value = 325
$ws.Range("K100").Value = 650
$ws.Range("M100").Value = -109

$ws.Range("H126").Value = 1277.1
$ws.Range("I126").Value = 1356
$ws.Range("J126").Value = 1198.2
$ws.Range("K126").Value = 4068
$ws.Range("L126").Value = 3594.6
$ws.Range("M126").Value = -1598
$ws.Range("N126").Value = -8534.6

$ws.Range("H132").Value = 6622
$ws.Range("I132").Value = 6607.1665
$ws.Range("J132").Value = 6666.5
$ws.Range("K132").Value = 19821.4995
$ws.Range("L132").Value = 19999.5
$ws.Range("M132").Value = -17291.4995
$ws.Range("N132").Value = -25059.5

$ws.Range("H136").Value = 4124.0244
$ws.Range("I136").Value = 1612.5
$ws.Range("J136").Value = 6089.5654
$ws.Range("K136").Value = 4837.5
$ws.Range("L136").Value = 18268.6962
$ws.Range("M136").Value = -2287.5
$ws.Range("N136").Value = -23368.6962
